# Add "American Journal of Transplantation" to the journal names map.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last data row (row 75) down into the new row 76 so the
# new row inherits the same cell formatting (style) the sheet already
# uses for this column, then overwrite its values with the new entry.
$ws.Rows(75).Copy()
$ws.Rows(76).Insert()

$ws.Range("A76").Value = "American Journal of Transplantation : Official Journal of the American Society of Transplantation and the American Society of Transplant Surgeons"
$ws.Range("B76").Value = "American Journal of Transplantation"

# Match the selection left behind by the edit.
$ws.Range("A77").Select()
